# Close and reopen browser between methods (test data update):
# Append the new "The Dead and the Missing" keyword rows to the HFHS_Keywords sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ASIN / Title / Keyword rows to append below the existing data
$newRows = @(
    @("B00XNZ0IZQ", "The Dead and the Missing", "Thrillers"),
    @("B00XNZ0IZQ", "The Dead and the Missing", "Private Detective"),
    @("B00XNZ0IZQ", "The Dead and the Missing", "Private Detective Series"),
    @("B00XNZ0IZQ", "The Dead and the Missing", "Crime Fiction")
)

$startRow = 4
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}

# Auto fit the ASIN/Title columns now that longer values were added
$ws.Range("A1:B7").EntireColumn.AutoFit()

# Leave the selection where the next write would go
$ws.Range("C10").Select()
